$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 394.0909  # H33 was 405.0909
$ws.Cells.Item(33, 9).Value = 266.875  # I33 was 282
$ws.Cells.Item(33, 11).Value = 266.875  # K33 was 282
$ws.Cells.Item(33, 13).Value = -37.875  # M33 was -53
$ws.Cells.Item(74, 8).Value = 15714  # H74 was 17199.6
$ws.Cells.Item(74, 9).Value = 12499.5  # I74 was 12999
$ws.Cells.Item(74, 11).Value = 12499.5  # K74 was 12999
$ws.Cells.Item(74, 13).Value = -11563.5  # M74 was -12063
$ws.Cells.Item(77, 8).Value = 15714  # H77 was 17199.6
$ws.Cells.Item(77, 9).Value = 12499.5  # I77 was 12999
$ws.Cells.Item(77, 11).Value = 62497.5  # K77 was 64995
$ws.Cells.Item(77, 13).Value = -57817.5  # M77 was -60315
$ws.Cells.Item(93, 8).Value = 46247  # H93 was 45996.4
$ws.Cells.Item(93, 10).Value = 46247  # J93 was 45996.4
$ws.Cells.Item(93, 12).Value = 46247  # L93 was 45996.4
$ws.Cells.Item(93, 14).Value = -51239  # N93 was -50988.4
$ws.Cells.Item(95, 8).Value = 62655.75  # H95 was 60124.4
$ws.Cells.Item(95, 10).Value = 62655.75  # J95 was 60124.4
$ws.Cells.Item(95, 12).Value = 62655.75  # L95 was 60124.4
$ws.Cells.Item(95, 14).Value = -68147.75  # N95 was -65616.39999999999
$ws.Cells.Item(106, 8).Value = 404.5  # H106 was 404
$ws.Cells.Item(106, 9).Value = 404.5  # I106 was 404
$ws.Cells.Item(106, 11).Value = 404.5  # K106 was 404
$ws.Cells.Item(106, 13).Value = 226.5  # M106 was 227
$ws.Cells.Item(115, 8).Value = 2086.75  # H115 was 2122.375
$ws.Cells.Item(115, 9).Value = 1115.6666  # I115 was 1163.1666
$ws.Cells.Item(115, 11).Value = 3346.9998  # K115 was 3489.4998
$ws.Cells.Item(115, 13).Value = -1779.9998  # M115 was -1922.4998
$ws.Cells.Item(130, 8).Value = 73567.2  # H130 was 66220.71000000001
$ws.Cells.Item(130, 10).Value = 125000  # J130 was 80000
$ws.Cells.Item(130, 12).Value = 125000  # L130 was 80000
$ws.Cells.Item(130, 14).Value = -135040  # N130 was -90040
$ws.Cells.Item(137, 8).Value = 2167.484  # H137 was 2152.0952
$ws.Cells.Item(137, 10).Value = 2280.2778  # J137 was 2223.3157
$ws.Cells.Item(137, 12).Value = 6840.8334  # L137 was 6669.9471
$ws.Cells.Item(137, 14).Value = -11940.8334  # N137 was -11769.9471
$ws.Cells.Item(138, 8).Value = 4435.88  # H138 was 4550.31
$ws.Cells.Item(138, 9).Value = 2655.9  # I138 was 2776.889
$ws.Cells.Item(138, 10).Value = 4633.656  # J138 was 4725.703
$ws.Cells.Item(138, 11).Value = 7967.700000000001  # K138 was 8330.667000000001
$ws.Cells.Item(138, 12).Value = 13900.968  # L138 was 14177.109
$ws.Cells.Item(138, 13).Value = -2827.700000000001  # M138 was -3190.667000000001
$ws.Cells.Item(138, 14).Value = -24180.968  # N138 was -24457.109
$ws.Cells.Item(141, 8).Value = 3026.818  # H141 was 3029.6
$ws.Cells.Item(141, 9).Value = 2498  # I141 was 2414.5
$ws.Cells.Item(141, 11).Value = 7494  # K141 was 7243.5
$ws.Cells.Item(141, 13).Value = -2314  # M141 was -2063.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14434.849  # H32 was 15240.693
$ws.Cells.Item(32, 9).Value = 6051.48  # I32 was 6365.1523
$ws.Cells.Item(32, 10).Value = 40632.875  # J32 was 40757.875
$ws.Cells.Item(32, 11).Value = 6051.48  # K32 was 6365.1523
$ws.Cells.Item(32, 12).Value = 40632.875  # L32 was 40757.875
$ws.Cells.Item(32, 13).Value = -5764.48  # M32 was -6078.1523
$ws.Cells.Item(32, 14).Value = -41206.875  # N32 was -41331.875
$ws.Cells.Item(94, 8).Value = 19998.666  # H94 was 19999
$ws.Cells.Item(94, 10).Value = 19998.666  # J94 was 19999
$ws.Cells.Item(94, 12).Value = 19998.666  # L94 was 19999
$ws.Cells.Item(94, 14).Value = -21800.666  # N94 was -21801
$ws.Cells.Item(132, 8).Value = 2695.8  # H132 was 1760.9828
$ws.Cells.Item(132, 9).Value = 2279.1  # I132 was 1436.9246
$ws.Cells.Item(132, 11).Value = 6837.299999999999  # K132 was 4310.7738
$ws.Cells.Item(132, 13).Value = -4307.299999999999  # M132 was -1780.7738

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 26388  # H76 was 25850
$ws.Cells.Item(76, 10).Value = 26388  # J76 was 25850
$ws.Cells.Item(76, 12).Value = 26388  # L76 was 25850
$ws.Cells.Item(76, 14).Value = -27018  # N76 was -26480
$ws.Cells.Item(79, 8).Value = 26388  # H79 was 25850
$ws.Cells.Item(79, 10).Value = 26388  # J79 was 25850
$ws.Cells.Item(79, 12).Value = 26388  # L79 was 25850
$ws.Cells.Item(79, 14).Value = -28572  # N79 was -28034
$ws.Cells.Item(94, 8).Value = 2699.8  # H94 was 2790.8
$ws.Cells.Item(94, 9).Value = 2877.5557  # I94 was 3385.5715
$ws.Cells.Item(94, 10).Value = 1100  # J94 was 1403
$ws.Cells.Item(94, 11).Value = 2877.5557  # K94 was 3385.5715
$ws.Cells.Item(94, 12).Value = 1100  # L94 was 1403
$ws.Cells.Item(94, 13).Value = -2426.5557  # M94 was -2934.5715
$ws.Cells.Item(94, 14).Value = -2002  # N94 was -2305
$ws.Cells.Item(96, 8).Value = 21320.2  # H96 was 20775.2
$ws.Cells.Item(96, 9).Value = 19865.223  # I96 was 19259.666
$ws.Cells.Item(96, 11).Value = 19865.223  # K96 was 19259.666
$ws.Cells.Item(96, 13).Value = -17119.223  # M96 was -16513.666
$ws.Cells.Item(97, 8).Value = 6359.5  # H97 was 8508.454
$ws.Cells.Item(97, 10).Value = 29999  # J97 was 29998.5
$ws.Cells.Item(97, 12).Value = 29999  # L97 was 29998.5
$ws.Cells.Item(97, 14).Value = -31981  # N97 was -31980.5
$ws.Cells.Item(99, 8).Value = 4207.241  # H99 was 4046.7742
$ws.Cells.Item(99, 9).Value = 2906.9375  # I99 was 2775.0557
$ws.Cells.Item(99, 11).Value = 2906.9375  # K99 was 2775.0557
$ws.Cells.Item(99, 13).Value = -1408.9375  # M99 was -1277.0557
$ws.Cells.Item(134, 8).Value = 1370.7307  # H134 was 1440.5834
$ws.Cells.Item(134, 9).Value = 1370.7307  # I134 was 1469.174
$ws.Cells.Item(134, 10).Value = 0  # J134 was 783
$ws.Cells.Item(134, 11).Value = 4112.1921  # K134 was 4407.522
$ws.Cells.Item(134, 12).Value = 0  # L134 was 2349
$ws.Cells.Item(134, 13).Value = -1577.1921  # M134 was -1872.522
$ws.Cells.Item(134, 14).ClearContents()  # N134 was -7419

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 4700.5  # H62 was 3705.2
$ws.Cells.Item(62, 9).Value = 2768  # I62 was 2500.5
$ws.Cells.Item(62, 10).Value = 6633  # J62 was 5512.25
$ws.Cells.Item(62, 11).Value = 2768  # K62 was 2500.5
$ws.Cells.Item(62, 12).Value = 6633  # L62 was 5512.25
$ws.Cells.Item(62, 13).Value = -2144  # M62 was -1876.5
$ws.Cells.Item(62, 14).Value = -7881  # N62 was -6760.25
$ws.Cells.Item(65, 8).Value = 4700.5  # H65 was 3705.2
$ws.Cells.Item(65, 9).Value = 2768  # I65 was 2500.5
$ws.Cells.Item(65, 10).Value = 6633  # J65 was 5512.25
$ws.Cells.Item(65, 11).Value = 13840  # K65 was 12502.5
$ws.Cells.Item(65, 12).Value = 33165  # L65 was 27561.25
$ws.Cells.Item(65, 13).Value = -10720  # M65 was -9382.5
$ws.Cells.Item(65, 14).Value = -39405  # N65 was -33801.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 8458.333000000001  # H80 was 7972.1113
$ws.Cells.Item(80, 9).Value = 6500  # I80 was 6749.5
$ws.Cells.Item(80, 10).Value = 8850  # J80 was 8321.429
$ws.Cells.Item(80, 11).Value = 19500  # K80 was 20248.5
$ws.Cells.Item(80, 12).Value = 26550  # L80 was 24964.287
$ws.Cells.Item(80, 13).Value = -18564  # M80 was -19312.5
$ws.Cells.Item(80, 14).Value = -28422  # N80 was -26836.287
$ws.Cells.Item(83, 8).Value = 8458.333000000001  # H83 was 7972.1113
$ws.Cells.Item(83, 9).Value = 6500  # I83 was 6749.5
$ws.Cells.Item(83, 10).Value = 8850  # J83 was 8321.429
$ws.Cells.Item(83, 11).Value = 58500  # K83 was 60745.5
$ws.Cells.Item(83, 12).Value = 79650  # L83 was 74892.861
$ws.Cells.Item(83, 13).Value = -53820  # M83 was -56065.5
$ws.Cells.Item(83, 14).Value = -89010  # N83 was -84252.861
$ws.Cells.Item(110, 8).Value = 20504  # H110 was 29666.334
$ws.Cells.Item(110, 9).Value = 20504  # I110 was 29666.334
$ws.Cells.Item(110, 11).Value = 61512  # K110 was 88999.00199999999
$ws.Cells.Item(110, 13).Value = -57422  # M110 was -84909.00199999999
$ws.Cells.Item(122, 8).Value = 2268  # H122 was 2448.125
$ws.Cells.Item(122, 9).Value = 2286.8462  # I122 was 2437.3
$ws.Cells.Item(122, 10).Value = 2251.6667  # J122 was 2455.8572
$ws.Cells.Item(122, 11).Value = 20581.6158  # K122 was 21935.7
$ws.Cells.Item(122, 12).Value = 20265.0003  # L122 was 22102.7148
$ws.Cells.Item(122, 13).Value = -18131.6158  # M122 was -19485.7
$ws.Cells.Item(122, 14).Value = -25165.0003  # N122 was -27002.7148
$ws.Cells.Item(129, 8).Value = 11000550  # H129 was 11000574
$ws.Cells.Item(129, 9).Value = 11000550  # I129 was 12375546
$ws.Cells.Item(129, 10).Value = 0  # J129 was 800
$ws.Cells.Item(129, 11).Value = 33001650  # K129 was 37126638
$ws.Cells.Item(129, 12).Value = 0  # L129 was 2400
$ws.Cells.Item(129, 13).Value = -32996650  # M129 was -37121638
$ws.Cells.Item(129, 14).ClearContents()  # N129 was -12400

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7748  # H80 was 6816.3335
$ws.Cells.Item(80, 9).Value = 4746.6665  # I80 was 5725
$ws.Cells.Item(80, 10).Value = 9999  # J80 was 8999
$ws.Cells.Item(80, 11).Value = 4746.6665  # K80 was 5725
$ws.Cells.Item(80, 12).Value = 9999  # L80 was 8999
$ws.Cells.Item(80, 13).Value = -3748.6665  # M80 was -4727
$ws.Cells.Item(80, 14).Value = -11995  # N80 was -10995
$ws.Cells.Item(83, 8).Value = 7748  # H83 was 6816.3335
$ws.Cells.Item(83, 9).Value = 4746.6665  # I83 was 5725
$ws.Cells.Item(83, 10).Value = 9999  # J83 was 8999
$ws.Cells.Item(83, 11).Value = 23733.3325  # K83 was 28625
$ws.Cells.Item(83, 12).Value = 49995  # L83 was 44995
$ws.Cells.Item(83, 13).Value = -18741.3325  # M83 was -23633
$ws.Cells.Item(83, 14).Value = -59979  # N83 was -54979
$ws.Cells.Item(98, 8).Value = 3642  # H98 was 3643
$ws.Cells.Item(98, 10).Value = 3642  # J98 was 3643
$ws.Cells.Item(98, 12).Value = 3642  # L98 was 3643
$ws.Cells.Item(98, 14).Value = -9632  # N98 was -9633
$ws.Cells.Item(100, 8).Value = 28036  # H100 was 28553.334
$ws.Cells.Item(100, 10).Value = 28036  # J100 was 28553.334
$ws.Cells.Item(100, 12).Value = 28036  # L100 was 28553.334
$ws.Cells.Item(100, 14).Value = -30200  # N100 was -30717.334
$ws.Cells.Item(104, 8).Value = 35671  # H104 was 0
$ws.Cells.Item(104, 10).Value = 35671  # J104 was 0
$ws.Cells.Item(104, 12).Value = 35671  # L104 was 0
$ws.Cells.Item(104, 14).Value = -42659  # N104 was None
$ws.Cells.Item(122, 8).Value = 2298.5417  # H122 was 2608.95
$ws.Cells.Item(122, 9).Value = 2011.7333  # I122 was 2289.9167
$ws.Cells.Item(122, 10).Value = 2776.5557  # J122 was 3087.5
$ws.Cells.Item(122, 11).Value = 6035.199900000001  # K122 was 6869.750100000001
$ws.Cells.Item(122, 12).Value = 8329.667099999999  # L122 was 9262.5
$ws.Cells.Item(122, 13).Value = -3585.199900000001  # M122 was -4419.750100000001
$ws.Cells.Item(122, 14).Value = -13229.6671  # N122 was -14162.5
$ws.Cells.Item(132, 8).Value = 4109.933  # H132 was 4995.8335
$ws.Cells.Item(132, 9).Value = 2093.625  # I132 was 3010
$ws.Cells.Item(132, 11).Value = 6280.875  # K132 was 9030
$ws.Cells.Item(132, 13).Value = -3750.875  # M132 was -6500

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2118.3333  # H61 was 2159.848
$ws.Cells.Item(61, 9).Value = 1513.0555  # I61 was 1621.325
$ws.Cells.Item(61, 11).Value = 1513.0555  # K61 was 1621.325
$ws.Cells.Item(61, 13).Value = -1311.0555  # M61 was -1419.325
$ws.Cells.Item(98, 8).Value = 25666.334  # H98 was 26666.334
$ws.Cells.Item(98, 9).Value = 21999  # I98 was 0
$ws.Cells.Item(98, 10).Value = 27500  # J98 was 26666.334
$ws.Cells.Item(98, 11).Value = 21999  # K98 was 0
$ws.Cells.Item(98, 12).Value = 27500  # L98 was 26666.334
$ws.Cells.Item(98, 13).Value = -19004  # M98 was None
$ws.Cells.Item(98, 14).Value = -33490  # N98 was -32656.334
$ws.Cells.Item(100, 8).Value = 668333.3  # H100 was 1001250
$ws.Cells.Item(100, 10).Value = 668333.3  # J100 was 1001250
$ws.Cells.Item(100, 12).Value = 668333.3  # L100 was 1001250
$ws.Cells.Item(100, 14).Value = -669415.3  # N100 was -1002332
$ws.Cells.Item(113, 8).Value = 2118.3333  # H113 was 2159.848
$ws.Cells.Item(113, 9).Value = 1513.0555  # I113 was 1621.325
$ws.Cells.Item(113, 11).Value = 1513.0555  # K113 was 1621.325
$ws.Cells.Item(113, 13).Value = 656.9445000000001  # M113 was 548.675
$ws.Cells.Item(132, 8).Value = 2656.9285  # H132 was 2728.8096
$ws.Cells.Item(132, 9).Value = 2138.9285  # I132 was 2246.75
$ws.Cells.Item(132, 11).Value = 6416.7855  # K132 was 6740.25
$ws.Cells.Item(132, 13).Value = -3886.7855  # M132 was -4210.25
$ws.Cells.Item(134, 8).Value = 86400  # H134 was 87000
$ws.Cells.Item(134, 10).Value = 86400  # J134 was 87000
$ws.Cells.Item(134, 12).Value = 86400  # L134 was 87000
$ws.Cells.Item(134, 14).Value = -96540  # N134 was -97140
$ws.Cells.Item(136, 8).Value = 4512.48  # H136 was 5034.2856
$ws.Cells.Item(136, 9).Value = 4032.4092  # I136 was 4534.5
$ws.Cells.Item(136, 11).Value = 12097.2276  # K136 was 13603.5
$ws.Cells.Item(136, 13).Value = -9547.2276  # M136 was -11053.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 50589.8  # H49 was 106031
$ws.Cells.Item(49, 10).Value = 50589.8  # J49 was 106031
$ws.Cells.Item(49, 12).Value = 50589.8  # L49 was 106031
$ws.Cells.Item(49, 14).Value = -51049.8  # N49 was -106491
$ws.Cells.Item(96, 8).Value = 1142.7142  # H96 was 1124.875
$ws.Cells.Item(96, 10).Value = 750  # J96 was 875
$ws.Cells.Item(96, 12).Value = 750  # L96 was 875
$ws.Cells.Item(96, 14).Value = -3496  # N96 was -3621
$ws.Cells.Item(100, 8).Value = 637.9  # H100 was 607.1818
$ws.Cells.Item(100, 9).Value = 396.66666  # I100 was 382.85715
$ws.Cells.Item(100, 11).Value = 793.33332  # K100 was 765.7143
$ws.Cells.Item(100, 13).Value = -252.33332  # M100 was -224.7143
$ws.Cells.Item(113, 8).Value = 286999.78  # H113 was 287002.66
$ws.Cells.Item(113, 9).Value = 1062.8518  # I113 was 1066.5927
$ws.Cells.Item(113, 11).Value = 3188.5554  # K113 was 3199.7781
$ws.Cells.Item(113, 13).Value = -1018.5554  # M113 was -1029.7781
$ws.Cells.Item(126, 8).Value = 2838.1538  # H126 was 2770.1853
$ws.Cells.Item(126, 9).Value = 2547.05  # I126 was 2473.524
$ws.Cells.Item(126, 11).Value = 7641.150000000001  # K126 was 7420.572
$ws.Cells.Item(126, 13).Value = -5171.150000000001  # M126 was -4950.572
$ws.Cells.Item(132, 8).Value = 2261.0667  # H132 was 2380.093
$ws.Cells.Item(132, 9).Value = 2056.8108  # I132 was 2191.3713
$ws.Cells.Item(132, 11).Value = 6170.432400000001  # K132 was 6574.113899999999
$ws.Cells.Item(132, 13).Value = -3640.432400000001  # M132 was -4044.113899999999
